# Fixed opening and closing shift bug, adding serialization and deserialization of objects
#
# The SSM sheet tracks which day(s) of the week (Mon=B .. Sun=H) apply to
# each opening/closing shift time-slot row. The opening-shift rows (7-14)
# had "Friday" (F) incorrectly flagged instead of "Tuesday" (C); the
# mid-day rows (15-18) were missing a day flag altogether (B/C/D should be
# 1); and several closing-shift rows (19-25) were missing the "Tuesday"
# (C) flag. This patch corrects those flags and also refreshes the
# active-cell selection left over from editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Opening-shift rows: Tuesday (C) should be on, Friday (F) should be off ---
foreach ($r in 7..14) {
    $ws.Cells.Item($r, 3).Value = 1   # column C = Tuesday
    $ws.Cells.Item($r, 6).Value = 0   # column F = Friday
}

# --- Mid-day rows: Monday/Tuesday/Wednesday (B/C/D) should be on ---
foreach ($r in 15..18) {
    $ws.Cells.Item($r, 2).Value = 1   # column B = Monday
    $ws.Cells.Item($r, 3).Value = 1   # column C = Tuesday
    $ws.Cells.Item($r, 4).Value = 1   # column D = Wednesday
}

# --- Closing-shift rows: Tuesday (C) should be on ---
foreach ($r in 19..24) {
    $ws.Cells.Item($r, 3).Value = 1   # column C = Tuesday
}

# --- Row 25: Monday/Tuesday/Wednesday (B/C/D) should be on ---
$ws.Cells.Item(25, 2).Value = 1   # column B = Monday
$ws.Cells.Item(25, 3).Value = 1   # column C = Tuesday
$ws.Cells.Item(25, 4).Value = 1   # column D = Wednesday

# --- Refresh the view: selection moves from D25 to D21, and the frozen/
#     scrolled top-left cell (previously A4) resets to the default (A1) ---
$ws.Range("A1").Select()
$ws.Range("D21").Select()
